$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- JIRA time-tracking data added for the three use-cases that previously
# --- had no hours logged yet (See Maps / Set Tags / Administrate Bars).
# Columns: C=Documentation, D=UI Design, E=Database, F=Warmup Phase,
#          G=Coding, H=Testing  (all stored as fraction-of-a-day durations)

# Row 13 - See Maps
$ws.Range("C13").Value = 60/1440
$ws.Range("D13").Value = 180/1440
$ws.Range("E13").Value = 60/1440
$ws.Range("F13").Value = 60/1440
$ws.Range("G13").Value = 60/1440
$ws.Range("H13").Value = 60/1440
$ws.Range("M13").Value = "The UI was very different hand had to be reworked serveral times to serve function as expected"

# Row 14 - Set Tags
$ws.Range("C14").Value = 15/1440
$ws.Range("D14").Value = 30/1440
$ws.Range("E14").Value = 10/1440
$ws.Range("F14").Value = 0/1440
$ws.Range("G14").Value = 64/1440
$ws.Range("H14").Value = 10/1440

# Row 15 - Administrate Bars
$ws.Range("C15").Value = 30/1440
$ws.Range("D15").Value = 80/1440
$ws.Range("E15").Value = 20/1440
$ws.Range("F15").Value = 0/1440
$ws.Range("G15").Value = 140/1440
$ws.Range("H15").Value = 10/1440

# --- Tabelle15 (Velocity helper table) now folds the two freshly-logged
# --- rows (14:15) into its "w/o outliers" totals, same as row 12 already was.
$ws.Range("O24").Formula = "=SUM(I3:I10,I12,I14:I15)*24"
$ws.Range("P24").Formula = "=SUM(L3:L10,L12,L14:L15)"

# --- Recalculate everything so dependent Estimate/FP formulas throughout
# --- the sheet (Tabelle1!K, and the FP-vs-Velocity helper rows) pick up
# --- the new Velocity.
$excel.CalculateFull()

# --- Leave the selection where the author ended up after the edit.
$ws.Range("Q24").Select() | Out-Null
